$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the outlier value (21) that was in A3 by deleting the entire row
# and shifting all subsequent rows up by one.
$ws.Rows.Item(3).Delete()
